# Applies the "Tratando de arreglar el formato" edit:
#  - shortens a couple of shared-string labels
#  - updates several computed values on the "reserva_total.prn" sheet
#  - adds bold/centered/bordered header formatting + fixed column widths on
#    the Pmax_Pgen / Mayor_maxima / Menor_optima sheets
#  - cleans up floating point noise in those same sheets
#  - merges A1:F1 and shrinks the A29 merge on reserva_total.prn

$wb = $excel.ActiveWorkbook

$sReserva = $wb.Worksheets.Item("reserva_total.prn")
$sPmax    = $wb.Worksheets.Item("Pmax_Pgen.prn")
$sMayor   = $wb.Worksheets.Item("Mayor_maxima.prn")
$sMenor   = $wb.Worksheets.Item("Menor_optima.prn")

# ---------------------------------------------------------------------------
# 1) Text tweaks (these rewrite the shared-string table entries)
# ---------------------------------------------------------------------------
$sReserva.Range("A29").Value = "LUEGO DEL RECORTE DE PONTECIA MAXIMA "

$sPmax.Range("B2").Value = "NUC-A"
$sPmax.Range("B3").Value = "NUC-B"
$sPmax.Range("B4").Value = "URBGEN"
$sPmax.Range("B5").Value = "HYDRO_G"

$sMayor.Range("B2").Value = "NUC-A"
$sMayor.Range("B3").Value = "NUC-B"
$sMayor.Range("B4").Value = "HYDRO_G"

$sMenor.Range("B2").Value = "URBGEN"

# ---------------------------------------------------------------------------
# 2) Updated figures in the second ("recorte") block of reserva_total.prn
# ---------------------------------------------------------------------------
$sReserva.Range("D31").Value = 195
$sReserva.Range("D32").Value = 385.68
$sReserva.Range("D33").Value = 580.6800000000001
$sReserva.Range("F34").Value = 17.83
$sReserva.Range("F40").Value = 157.5
$sReserva.Range("F41").Value = 4.84
$sReserva.Range("D43").Value = 945
$sReserva.Range("D44").Value = 1615.67
$sReserva.Range("D46").Value = 3460.67
$sReserva.Range("D49").Value = 285.67
$sReserva.Range("D51").Value = 385.67
$sReserva.Range("F53").Value = 480.68

# ---------------------------------------------------------------------------
# 3) Merge / layout tweaks on reserva_total.prn
# ---------------------------------------------------------------------------
$sReserva.Range("A1:F1").Merge()
$sReserva.Range("A29:H29").UnMerge()
$sReserva.Range("A29:F29").Merge()

$sReserva.Range("A1").HorizontalAlignment = -4108
$sReserva.Range("A29").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 4) Header / data formatting for the three generator-list sheets
# ---------------------------------------------------------------------------
foreach ($ws in @($sPmax, $sMayor, $sMenor)) {
    $used = $ws.Range("A1:I1")
    $headerRow = $ws.Range("A1:I1")
    $headerRow.Font.Bold = $true
    $headerRow.HorizontalAlignment = -4108
    $headerRow.Borders.LineStyle = 1

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $dataRange = $ws.Range("A2:I" + $lastRow)
    $dataRange.HorizontalAlignment = -4108
    $dataRange.Borders.LineStyle = 1

    for ($col = 1; $col -le 9; $col++) {
        $ws.Columns.Item($col).ColumnWidth = 20
    }
}

# ---------------------------------------------------------------------------
# 5) Clean up floating point noise on Pmax_Pgen.prn
# ---------------------------------------------------------------------------
$sPmax.Range("D2").Value = 945
$sPmax.Range("E2").Value = 750
$sPmax.Range("F2").Value = 195
$sPmax.Range("D3").Value = 945
$sPmax.Range("E3").Value = 750
$sPmax.Range("F3").Value = 195
$sPmax.Range("D4").Value = 900
$sPmax.Range("E4").Value = 800
$sPmax.Range("F4").Value = 100
$sPmax.Range("D5").Value = 800.4
$sPmax.Range("E5").Value = 580
$sPmax.Range("F5").Value = 220.4

# ---------------------------------------------------------------------------
# 6) Clean up floating point noise on Mayor_maxima.prn
# ---------------------------------------------------------------------------
$sMayor.Range("D2").Value = 945
$sMayor.Range("E2").Value = 750
$sMayor.Range("F2").Value = 195
$sMayor.Range("D3").Value = 945
$sMayor.Range("E3").Value = 750
$sMayor.Range("F3").Value = 195
$sMayor.Range("D4").Value = 800.4
$sMayor.Range("E4").Value = 580
$sMayor.Range("F4").Value = 220.4

# ---------------------------------------------------------------------------
# 7) Clean up floating point noise on Menor_optima.prn
# ---------------------------------------------------------------------------
$sMenor.Range("D2").Value = 900
$sMenor.Range("E2").Value = 800
$sMenor.Range("F2").Value = 100

Write-Output "edit applied"
